$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.637.38'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '1.632.26'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D5").Value = '212.74'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("E6").Value = '  +3.17%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '0.252'
$ws.Range("E8").Value = '  +1.91%  '
$ws.Range("D9").Value = '0.0621'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = '19.13'
$ws.Range("E10").Value = '  +2.12%  '
$ws.Range("D11").Value = '0.0843'
$ws.Range("D12").Value = '1.860.54'
$ws.Range("D13").Value = '1.611.40'
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").Value = '0.523'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").Value = '26.637.45'
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("D17").Value = '63.15'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("D19").Value = '217.54'
$ws.Range("E19").Value = '  +7.91%  '
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("D23").Value = '9.37'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D25").Value = '147.92'
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("D28").Value = '6.86'
$ws.Range("E28").Value = '  +4.67%  '
$ws.Range("D29").Value = '15.45'
$ws.Range("E29").Value = '  +1.88%  '
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("E32").Value = '  +3.77%  '
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").Value = '1.208.43'
$ws.Range("E36").Value = '  +2.68%  '
$ws.Range("E37").Value = '  +5.46%  '
$ws.Range("D38").Value = '0.806'
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("E41").Value = '  -1.96%  '
$ws.Range("D42").Value = '5.40'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").Value = '1.772.28'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("D45").Value = '92.46'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").Value = '54.62'
$ws.Range("E47").Value = '  +1.64%  '
$ws.Range("D48").Value = '0.0513'
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("D49").Value = '7.61'
$ws.Range("E49").Value = '  +4.54%  '
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("E51").Value = '  +0.16%  '
